$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected ICDC "FilesTab" Cypher query (row 4 / cell B4):
# - removed the "File Type" projection line
# - removed the "Breed" projection line
$newQuery = @"
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Irish Wolfhound'] 
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS ``File Name``,
         coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
        coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_size, '') AS ``Size``,
        coalesce(c.case_id, '') AS ``Case ID``,
         coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS ``Study Code``
"@

$ws.Range("B4").Value = $newQuery

# The query text got shorter, so the wrapped row height shrinks accordingly.
$ws.Rows.Item(4).RowHeight = 217.5

# Selection moved from C4 to B4.
$null = $ws.Range("B4").Select()
